$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "SR012-Humificador" with "WR005-Plastilina" (A2, P2)
$ws.Range("A2").Value = "WR005-Plastilina"
$ws.Range("P2").Value = "WR005-Plastilina"

# Replace "SB034-Porta" with "WX001-Juguete" (A3, P3)
$ws.Range("A3").Value = "WX001-Juguete"
$ws.Range("P3").Value = "WX001-Juguete"

# Replace "SR012" with "WR005" (B2, T2)
$ws.Range("B2").Value = "WR005"
$ws.Range("T2").Value = "WR005"

# Replace "SB034" with "WX001" (B3, T3)
$ws.Range("B3").Value = "WX001"
$ws.Range("T3").Value = "WX001"

# Update unit sale price for row 3 (G3) from 2 to 9.5
$ws.Range("G3").Value = 9.5
